# Auto-generated Excel COM-interop script
# Applies numeric updates to the Leve profit-tracking tables (per-sheet)
# as captured from the commit's scheduled Kraken_Profits price-refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 491.33334
$ws.Range("I2").Value = 491.33334
$ws.Range("K2").Value = 491.33334
$ws.Range("M2").Value = -378.33334
$ws.Range("H40").Value = 6790.1
$ws.Range("I40").Value = 3580.2
$ws.Range("K40").Value = 3580.2
$ws.Range("M40").Value = -3405.2
$ws.Range("H43").Value = 1999
$ws.Range("I43").Value = 1999
$ws.Range("K43").Value = 1999
$ws.Range("M43").Value = -1930
$ws.Range("H98").Value = 553
$ws.Range("I98").Value = 553
$ws.Range("K98").Value = 553
$ws.Range("M98").Value = 945
$ws.Range("H122").Value = 553
$ws.Range("I122").Value = 553
$ws.Range("K122").Value = 1659
$ws.Range("M122").Value = 791
$ws.Range("H127").Value = 918.8570999999999
$ws.Range("I127").Value = 872
$ws.Range("K127").Value = 2616
$ws.Range("M127").Value = 2344

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 162.75
$ws.Range("I5").Value = 162.75
$ws.Range("K5").Value = 162.75
$ws.Range("M5").Value = -50.75
$ws.Range("H61").Value = 2823.5
$ws.Range("I61").Value = 2615
$ws.Range("K61").Value = 2615
$ws.Range("M61").Value = -2403
$ws.Range("H74").Value = 3339.4443
$ws.Range("I74").Value = 3339.4443
$ws.Range("K74").Value = 3339.4443
$ws.Range("M74").Value = -2465.4443
$ws.Range("H77").Value = 3339.4443
$ws.Range("I77").Value = 3339.4443
$ws.Range("K77").Value = 16697.2215
$ws.Range("M77").Value = -12329.2215
$ws.Range("H92").Value = 39998.5
$ws.Range("J92").Value = 39998.5
$ws.Range("L92").Value = 39998.5
$ws.Range("N92").Value = -44990.5
$ws.Range("H110").Value = 2011.2
$ws.Range("I110").Value = 1678.1
$ws.Range("J110").Value = 2677.4
$ws.Range("K110").Value = 1678.1
$ws.Range("L110").Value = 2677.4
$ws.Range("M110").Value = 366.9000000000001
$ws.Range("N110").Value = -6767.4
$ws.Range("H136").Value = 2823.5
$ws.Range("I136").Value = 2615
$ws.Range("K136").Value = 7845
$ws.Range("M136").Value = -5295

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 162.75
$ws.Range("I4").Value = 162.75
$ws.Range("K4").Value = 162.75
$ws.Range("M4").Value = -47.75
$ws.Range("H15").Value = 45000
$ws.Range("I15").Value = 45000
$ws.Range("K15").Value = 45000
$ws.Range("M15").Value = -44773
$ws.Range("H107").Value = 11237.083
$ws.Range("I107").Value = 4531.4546
$ws.Range("J107").Value = 84999
$ws.Range("K107").Value = 4531.4546
$ws.Range("L107").Value = 84999
$ws.Range("M107").Value = -2611.4546
$ws.Range("N107").Value = -88839

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5000000
$ws.Range("I4").Value = 5000000
$ws.Range("K4").Value = 5000000
$ws.Range("M4").Value = -4999888
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502
$ws.Range("H105").Value = 920
$ws.Range("I105").Value = 920
$ws.Range("K105").Value = 920
$ws.Range("M105").Value = 827
$ws.Range("H122").Value = 1568.5
$ws.Range("I122").Value = 1568.5
$ws.Range("K122").Value = 4705.5
$ws.Range("M122").Value = -2255.5
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530
$ws.Range("H132").Value = 3926.111
$ws.Range("I132").Value = 3833.5715
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 11500.7145
$ws.Range("L132").Value = 12750
$ws.Range("M132").Value = -8970.7145
$ws.Range("N132").Value = -17810

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2862.2144
$ws.Range("I5").Value = 2718
$ws.Range("J5").Value = 2919.9
$ws.Range("K5").Value = 8154
$ws.Range("L5").Value = 8759.700000000001
$ws.Range("M5").Value = -8042
$ws.Range("N5").Value = -8983.700000000001
$ws.Range("H80").Value = 5788.3335
$ws.Range("J80").Value = 6087.273
$ws.Range("L80").Value = 18261.819
$ws.Range("N80").Value = -20133.819
$ws.Range("H81").Value = 5000
$ws.Range("I81").Value = 5000
$ws.Range("K81").Value = 15000
$ws.Range("M81").Value = -13877
$ws.Range("H83").Value = 5788.3335
$ws.Range("J83").Value = 6087.273
$ws.Range("L83").Value = 54785.457
$ws.Range("N83").Value = -64145.457
$ws.Range("H84").Value = 5000
$ws.Range("I84").Value = 5000
$ws.Range("K84").Value = 45000
$ws.Range("M84").Value = -39384
$ws.Range("H135").Value = 2862.2144
$ws.Range("I135").Value = 2718
$ws.Range("J135").Value = 2919.9
$ws.Range("K135").Value = 24462
$ws.Range("L135").Value = 26279.1
$ws.Range("M135").Value = -21927
$ws.Range("N135").Value = -31349.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H126").Value = 3285.25
$ws.Range("I126").Value = 2713.6667
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 8141.000100000001
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -5671.000100000001
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1798.3334
$ws.Range("I16").Value = 1798.3334
$ws.Range("K16").Value = 1798.3334
$ws.Range("M16").Value = -1628.3334
$ws.Range("H61").Value = 3098.1667
$ws.Range("I61").Value = 3098.1667
$ws.Range("K61").Value = 3098.1667
$ws.Range("M61").Value = -2896.1667
$ws.Range("H100").Value = 4484
$ws.Range("I100").Value = 4666.6665
$ws.Range("J100").Value = 4210
$ws.Range("K100").Value = 4666.6665
$ws.Range("L100").Value = 4210
$ws.Range("M100").Value = -4125.6665
$ws.Range("N100").Value = -5292
$ws.Range("H113").Value = 3098.1667
$ws.Range("I113").Value = 3098.1667
$ws.Range("K113").Value = 3098.1667
$ws.Range("M113").Value = -928.1667000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 5010000
$ws.Range("J5").Value = 5010000
$ws.Range("L5").Value = 5010000
$ws.Range("N5").Value = -5010224
$ws.Range("H107").Value = 625.0769
$ws.Range("I107").Value = 532.7
$ws.Range("K107").Value = 1598.1
$ws.Range("M107").Value = 321.8999999999999
$ws.Range("H132").Value = 3298.625
$ws.Range("I132").Value = 3657.7144
$ws.Range("K132").Value = 10973.1432
$ws.Range("M132").Value = -8443.143199999999

